# "semana 52 de 2025"
#
# The sheet is a weekly case-count matrix: row 1 holds the epidemiological
# week numbers (1..50 so far) as text headers, and each subsequent row is a
# health facility with one case-count value per week. This adds the two new
# weekly columns, week 51 (column BB) and week 52 (column BC), with the
# header labels and each facility's counts for those two weeks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new week-number columns ------------------------------
# The existing week headers (D1:BA1, "1".."50") are stored as text, not
# numbers, even though they look numeric. A leading apostrophe forces the
# same text storage here. Setting a numeric-looking value this way makes
# Excel tag the cell with a "stored as text" quote-prefix flag (a new
# style), so style is reset to Normal and bold+center is re-applied
# afterwards to land back on the same header style used by every other
# cell in row 1.
$ws.Range("BB1").Value = "'51"
$ws.Range("BB1").Style = "Normal"
$ws.Range("BB1").Font.Bold = $true
$ws.Range("BB1").HorizontalAlignment = -4108   # xlCenter

$ws.Range("BC1").Value = "'52"
$ws.Range("BC1").Style = "Normal"
$ws.Range("BC1").Font.Bold = $true
$ws.Range("BC1").HorizontalAlignment = -4108   # xlCenter

# --- Data rows: case counts for week 51 (BB) and week 52 (BC) ---------
# Only the rows/columns actually present in the source update are listed;
# rows not touched by the update (e.g. row 4, 12, 13...) are left as-is.
$weekData = @{
    2  = @{ BB = 0; BC = 0 };
    3  = @{ BB = 0; BC = 0 };
    5  = @{ BB = 0; BC = 0 };
    6  = @{ BB = 2; BC = 0 };
    7  = @{ BB = 0; BC = 0 };
    8  = @{ BB = 0; BC = 0 };
    9  = @{ BB = 0; BC = 0 };
    10 = @{ BB = 0; BC = 0 };
    11 = @{ BB = 0 };
    14 = @{ BB = 0; BC = 0 };
    15 = @{ BB = 0 };
    16 = @{ BB = 0; BC = 0 };
    17 = @{ BB = 0 };
    23 = @{ BB = 0 };
    25 = @{ BB = 0; BC = 0 };
    28 = @{ BB = 0; BC = 1 };
    29 = @{ BB = 0; BC = 0 };
    30 = @{ BB = 5; BC = 5 };
    31 = @{ BB = 0; BC = 0 };
    35 = @{ BB = 1; BC = 1 };
    36 = @{ BB = 1; BC = 1 };
    37 = @{ BB = 0; BC = 0 };
    38 = @{ BB = 0; BC = 0 };
    41 = @{ BB = 0; BC = 0 };
    42 = @{ BB = 0; BC = 0 };
    43 = @{ BB = 0; BC = 0 };
    45 = @{ BB = 0; BC = 0 };
    46 = @{ BB = 0; BC = 0 };
    47 = @{ BB = 0; BC = 0 };
    48 = @{ BB = 0; BC = 0 };
    49 = @{ BB = 0; BC = 0 };
    50 = @{ BB = 0; BC = 0 };
    51 = @{ BB = 0; BC = 0 };
    54 = @{ BB = 0; BC = 0 };
    55 = @{ BB = 0; BC = 0 };
    56 = @{ BB = 2; BC = 0 };
    57 = @{ BB = 0; BC = 0 };
    58 = @{ BB = 0; BC = 0 };
    59 = @{ BB = 0; BC = 0 };
}

foreach ($row in $weekData.Keys) {
    $cells = $weekData[$row]
    foreach ($col in $cells.Keys) {
        $ws.Range("$col$row").Value = $cells[$col]
    }
}
